$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix casing of the metadata4Ing header labels -> metadata4ing
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"

# Add the new RXNO_DEF column (header + data), matching the header style
# used by the other header cells (B1:E1).
$ws.Range("F1").Value = "RXNO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F2").Value = "['p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]', locstr(`"Process, i.e., a physical entity with a temporal evolution that 'has a meaning for the ontologist'`", 'en')]"
